$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G to match column F's width (16.21875 characters)
$ws.Columns("G").ColumnWidth = 15.333333333333334

# Row 8
$ws.Range("G8").Value2 = 0

# Row 9
$ws.Range("E9").Value2 = 0.90000000000000002
$ws.Range("G9").Value2 = 0

# Row 13
$ws.Range("F13").Value2 = 0.90000000000000002

# Row 14
$ws.Range("G14").Value2 = 0

# Row 17
$ws.Range("E17").Value2 = -3059.9540000000002
$ws.Range("F17").Value2 = -3059.9540000000002
$ws.Range("G17").Value2 = -109265021298.28999

# Row 18
$ws.Range("E18").Value2 = -2863.9119999999998
$ws.Range("F18").Value2 = -2863.9119999999998
$ws.Range("G18").Value2 = -105945556394.51401

# Row 19
$ws.Range("E19").Value2 = -430.22500000000002
$ws.Range("F19").Value2 = -430.22500000000002
$ws.Range("G19").Value2 = 215210571333.80499
$ws.Range("I19").Value2 = 1064227.1839999999

# Row 20
$ws.Range("E20").Value2 = -6354.0919999999996
$ws.Range("F20").Value2 = -6354.0919999999996
$ws.Range("G20").Value2 = -6358.9989999999998
$ws.Range("I20").Value2 = -23876.825000000001
